$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.327.85"
$ws.Range("E2").Value = "  -3.85%  "

# Row 3
$ws.Range("D3").Value = "2.462.29"
$ws.Range("E3").Value = "  -6.41%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.33"
$ws.Range("E5").Value = "  -5.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.43"
$ws.Range("E6").Value = "  -6.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("E8").Value = "  -6.74%  "

# Row 9
$ws.Range("D9").Value = "2.460.26"
$ws.Range("E9").Value = "  -6.45%  "

# Row 10
$ws.Range("E10").Value = "  -9.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.45"
$ws.Range("E11").Value = "  -6.32%  "

# Row 12
$ws.Range("E12").Value = "  -1.77%  "

# Row 13
$ws.Range("E13").Value = "  -8.44%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.12"
$ws.Range("E14").Value = "  -8.93%  "

# Row 15
$ws.Range("D15").Value = "2.905.44"
$ws.Range("E15").Value = "  -6.52%  "

# Row 16
$ws.Range("E16").Value = "  -9.61%  "

# Row 17
$ws.Range("D17").Value = "61.246.55"
$ws.Range("E17").Value = "  -3.86%  "

# Row 18
$ws.Range("D18").Value = "2.453.83"
$ws.Range("E18").Value = "  -7.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("E19").Value = "  -8.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.07"
$ws.Range("E20").Value = "  -8.34%  "

# Row 21
$ws.Range("E21").Value = "  -7.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "318.78"
$ws.Range("E22").Value = "  -7.21%  "

# Row 23
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.84"
$ws.Range("E24").Value = "  -3.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.87"
$ws.Range("E25").Value = "  -6.53%  "

# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.578.25"
$ws.Range("E26").Value = "  -6.88%  "

# Row 27
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "550.01"
$ws.Range("E27").Value = "  -5.38%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").Value = "  -14.40%  "

# Row 29
$ws.Range("E29").Value = "  -0.87%  "

# Row 30
$ws.Range("E30").Value = "  -11.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.25"
$ws.Range("E31").Value = "  -10.55%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.62"
$ws.Range("E32").Value = "  -7.97%  "

# Row 33
$ws.Range("E33").Value = "  -8.19%  "

# Row 34
$ws.Range("E34").Value = "  -8.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  -8.64%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.86"
$ws.Range("E36").Value = "  -11.87%  "

# Row 37
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.83"
$ws.Range("E38").Value = "  -11.70%  "

# Row 39
$ws.Range("E39").Value = "  -6.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.43"
$ws.Range("E40").Value = "  -6.82%  "

# Row 41
$ws.Range("E41").Value = "  -7.15%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.03"
$ws.Range("E42").Value = "  -8.61%  "

# Row 43
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.43"
$ws.Range("E44").Value = "  -3.87%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("E45").Value = "  -10.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.63"
$ws.Range("E46").Value = "  -9.98%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.60"
$ws.Range("E47").Value = "  -8.32%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.41"
$ws.Range("E48").Value = "  -10.95%  "

# Row 49
$ws.Range("E49").Value = "  -8.82%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.588"
$ws.Range("E50").Value = "  -7.38%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0937"
$ws.Range("E51").Value = "  -6.55%  "

